# The sheet "Hortaliza, Femacal de La Calera - Apio" gets one new weekly
# price-report row inserted before the existing row 391 (a new Primera-grade
# "Apio" entry dated 2022-07-27 / serial 44769), which pushes every
# subsequent row (old 391..430) down by one, growing the used range from
# A1:R430 to A1:R431. Old row 430's data simply lands on new row 431
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 391..430 down to 392..431, opening up a blank row 391.
$ws.Rows.Item(391).Insert()

# Populate the newly-opened row 391 with the new weekly record.
$ws.Range("A391").Value = 3
$ws.Range("B391").Value = "Femacal de La Calera"
$ws.Range("C391").Value = "Coquimbo"
$ws.Range("D391").Value = 44769
$ws.Range("E391").Value = 5
$ws.Range("F391").Value = 100112017
$ws.Range("G391").Value = "Apio"
$ws.Range("H391").Value = "Americana (o)"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 230
$ws.Range("K391").Value = 9000
$ws.Range("L391").Value = 9500
$ws.Range("M391").Value = 9261
$ws.Range("N391").Value = "$/docena de matas"
$ws.Range("O391").Value = "Pan de Azúcar"
$ws.Range("P391").Value = 1544
$ws.Range("Q391").Value = 6
$ws.Range("R391").Value = "Hortaliza"
